# Insert a new item card row for "时间跳跃" (Jet Lag) at row 19 of Sheet1,
# pushing the existing rows 19-44 down to 20-45.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 19 - this shifts rows 19:44 down to 20:45.
$ws.Rows.Item(19).Insert()

# Populate the new row 19 with the new card's data.
$ws.Cells.Item(19, 1).Value = "时间跳跃"
$ws.Cells.Item(19, 2).Value = $ws.Cells.Item(20, 2).Value2
$ws.Cells.Item(19, 3).Value = 2
$ws.Cells.Item(19, 4).Value = "牌通过移动后可以发动：将那张牌复位。"

$ws.Rows.Item(19).RowHeight = 71.25
